# إضافة حدث جديد في Card22
# Fills the previously-blank placeholder cells on row 15 with "nan",
# then appends a brand-new service event as row 16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card22")

# --- Row 15: the existing blank placeholder cells get the literal text "nan" ---
$ws.Range("B15").Value = "nan"
$ws.Range("C15").Value = "nan"
$ws.Range("D15").Value = "nan"
$ws.Range("E15").Value = "nan"
$ws.Range("F15").Value = "nan"
$ws.Range("G15").Value = "nan"
$ws.Range("H15").Value = "nan"
$ws.Range("I15").Value = "nan"
$ws.Range("J15").Value = "nan"
$ws.Range("K15").Value = "nan"
$ws.Range("M15").Value = "nan"

# --- Row 16: brand-new event entry for Card22 ---
# Force column A to be stored as text ("22"), matching the rest of the column.
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "22"

$ws.Range("L16").Value = "15\12\2024"
$ws.Range("M16").Value = "4320 h"
$ws.Range("N16").Value = "تم عمل صيانه نصف سنويه"
$ws.Range("O16").Value = "تيم العمل"
